$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the I/J columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy H1's formatting (bold font, border, centered alignment) onto the
# new header cells so they share the same cell style as the rest of row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# New data values for rows 2-4
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 4
